# da sua them giangvien bang file excel
# Rename the 4 faculty-code rows (A2:A5) from the old placeholder names
# (CNN / CATOON / HEHE / HAHA) to the new ones (Test1 / Test2 / Test3 / Test4),
# move the active selection to B8, and size up columns A and B so the new,
# slightly longer labels are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test1"
$ws.Range("A3").Value = "Test2"
$ws.Range("A4").Value = "Test3"
$ws.Range("A5").Value = "Test4"

$ws.Columns.Item(1).ColumnWidth = 18.44140625
$ws.Columns.Item(2).ColumnWidth = 15.109375

$ws.Range("B8").Select()
